# Extend the daily data table through 28/06/2021 (serial date 44375).
# The sheet currently ends at row 269 (date serial 44343, i.e. 27/05/2021).
# We append rows 270..301 with consecutive date serials, zeros in B/C/D,
# and the same date style (s="2") used by the existing date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 270
$endRow = 301
$startDate = 44344

for ($i = 0; $i -le ($endRow - $startRow); $i++) {
    $row = $startRow + $i
    $date = $startDate + $i

    # Copy the previous row's date cell first so the new cell inherits
    # its style (border/bold/date number format), then overwrite the value.
    $ws.Cells.Item($row - 1, 1).Copy($ws.Cells.Item($row, 1))
    $ws.Cells.Item($row, 1).Value = $date

    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
